# Update odds values in the "Jogos da Semana" worksheet
# per the scraped FlashScore refresh (2024-10-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Brusque vs Ituano)
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63

# Row 3 (Guarani vs CRB)
$ws.Range("G3").Value = 2.05
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 4
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 2.25
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 8.5
$ws.Range("AH3").Value = 19
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 51
$ws.Range("AP3").Value = 29
$ws.Range("AT3").Value = 2.25
$ws.Range("AV3").Value = 81
$ws.Range("AY3").Value = 41

# Row 4 (Patriotas vs Santa Fe)
$ws.Range("G4").Value = 3.1
$ws.Range("I4").Value = 2.35
$ws.Range("L4").Value = 3.2
$ws.Range("AK4").Value = 23
$ws.Range("AX4").Value = 15
